$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 1111, pushing existing rows 1111-1177 down to 1114-1180
$ws.Range("A1111:T1113").EntireRow.Insert()

# Fill in the 3 new rows with the new weekly data (Frutilla, Vega Central Mapocho de Santiago,
# Provincia de San Antonio, fecha 2023-01-05 / serial 44931)

# Row 1111: Calidad = Especial
$ws.Cells.Item(1111, 1).Value = 9
$ws.Cells.Item(1111, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1111, 3).Value = "Metropolitana"
$ws.Cells.Item(1111, 4).Value = 44931
$ws.Cells.Item(1111, 5).Value = 13
$ws.Cells.Item(1111, 6).Value = "Fruta"
$ws.Cells.Item(1111, 7).Value = 100101
$ws.Cells.Item(1111, 8).Value = "Berries"
$ws.Cells.Item(1111, 9).Value = 100112025
$ws.Cells.Item(1111, 10).Value = "Frutilla"
$ws.Cells.Item(1111, 11).Value = "Sin especificar"
$ws.Cells.Item(1111, 12).Value = "Especial"
$ws.Cells.Item(1111, 13).Value = 630
$ws.Cells.Item(1111, 14).Value = 6500
$ws.Cells.Item(1111, 15).Value = 7000
$ws.Cells.Item(1111, 16).Value = 6722
$ws.Cells.Item(1111, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(1111, 18).Value = "Provincia de San Antonio"
$ws.Cells.Item(1111, 19).Value = 960
$ws.Cells.Item(1111, 20).Value = 7

# Row 1112: Calidad = Primera
$ws.Cells.Item(1112, 1).Value = 9
$ws.Cells.Item(1112, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1112, 3).Value = "Metropolitana"
$ws.Cells.Item(1112, 4).Value = 44931
$ws.Cells.Item(1112, 5).Value = 13
$ws.Cells.Item(1112, 6).Value = "Fruta"
$ws.Cells.Item(1112, 7).Value = 100101
$ws.Cells.Item(1112, 8).Value = "Berries"
$ws.Cells.Item(1112, 9).Value = 100112025
$ws.Cells.Item(1112, 10).Value = "Frutilla"
$ws.Cells.Item(1112, 11).Value = "Sin especificar"
$ws.Cells.Item(1112, 12).Value = "Primera"
$ws.Cells.Item(1112, 13).Value = 600
$ws.Cells.Item(1112, 14).Value = 4500
$ws.Cells.Item(1112, 15).Value = 5000
$ws.Cells.Item(1112, 16).Value = 4750
$ws.Cells.Item(1112, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(1112, 18).Value = "Provincia de San Antonio"
$ws.Cells.Item(1112, 19).Value = 679
$ws.Cells.Item(1112, 20).Value = 7

# Row 1113: Calidad = Segunda
$ws.Cells.Item(1113, 1).Value = 9
$ws.Cells.Item(1113, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(1113, 3).Value = "Metropolitana"
$ws.Cells.Item(1113, 4).Value = 44931
$ws.Cells.Item(1113, 5).Value = 13
$ws.Cells.Item(1113, 6).Value = "Fruta"
$ws.Cells.Item(1113, 7).Value = 100101
$ws.Cells.Item(1113, 8).Value = "Berries"
$ws.Cells.Item(1113, 9).Value = 100112025
$ws.Cells.Item(1113, 10).Value = "Frutilla"
$ws.Cells.Item(1113, 11).Value = "Sin especificar"
$ws.Cells.Item(1113, 12).Value = "Segunda"
$ws.Cells.Item(1113, 13).Value = 500
$ws.Cells.Item(1113, 14).Value = 3000
$ws.Cells.Item(1113, 15).Value = 3500
$ws.Cells.Item(1113, 16).Value = 3220
$ws.Cells.Item(1113, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(1113, 18).Value = "Provincia de San Antonio"
$ws.Cells.Item(1113, 19).Value = 460
$ws.Cells.Item(1113, 20).Value = 7
